$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 526
$ws.Cells.Item(526, 1).Value = "21CRB01268"
$ws.Cells.Item(526, 2).Value = "Bunner"
$ws.Cells.Item(526, 3).Value = "POSSESSION DRUG PARAPHERNALIA"
$ws.Cells.Item(526, 4).Value = "2925.14(C)"
$ws.Cells.Item(526, 5).Value = "M4"
$ws.Cells.Item(526, 6).Value = "Guilty"
$ws.Cells.Item(526, 7).Value = "Guilty"
$c_526_H = $ws.Cells.Item(526, 8)
$c_526_H.NumberFormat = "@"
$c_526_H.Value = "$ 50"
$c_526_H.Style = "Normal"
$c_526_I = $ws.Cells.Item(526, 9)
$c_526_I.NumberFormat = "@"
$c_526_I.Value = "$ 25"
$c_526_I.Style = "Normal"

# Row 527
$ws.Cells.Item(527, 1).Value = "21CRB01268"
$ws.Cells.Item(527, 2).Value = "Bunner"
$ws.Cells.Item(527, 3).Value = "POSSESSION DRUG PARAPHERNALIA"
$ws.Cells.Item(527, 4).Value = "2925.14(C)"
$ws.Cells.Item(527, 5).Value = "M4"
$ws.Cells.Item(527, 6).Value = "Guilty"
$ws.Cells.Item(527, 7).Value = "Guilty"
$c_527_H = $ws.Cells.Item(527, 8)
$c_527_H.NumberFormat = "@"
$c_527_H.Value = "$ 0"
$c_527_H.Style = "Normal"
$c_527_I = $ws.Cells.Item(527, 9)
$c_527_I.NumberFormat = "@"
$c_527_I.Value = "$ 0"
$c_527_I.Style = "Normal"

# Row 528
$ws.Cells.Item(528, 1).Value = "21CRB01291"
$ws.Cells.Item(528, 2).Value = "Bunner"
$ws.Cells.Item(528, 3).Value = "PERMISSION REQ'D TO USE LICENSED DOCK"
$ws.Cells.Item(528, 4).Value = "1501:46-12-04"
$ws.Cells.Item(528, 5).Value = "MM"
$ws.Cells.Item(528, 6).Value = "Guilty"
$ws.Cells.Item(528, 7).Value = "Guilty"
$c_528_H = $ws.Cells.Item(528, 8)
$c_528_H.NumberFormat = "@"
$c_528_H.Value = "$ 0"
$c_528_H.Style = "Normal"
$c_528_I = $ws.Cells.Item(528, 9)
$c_528_I.NumberFormat = "@"
$c_528_I.Value = "$ 0"
$c_528_I.Style = "Normal"

# Row 529
$ws.Cells.Item(529, 1).Value = "21TRC08418"
$ws.Cells.Item(529, 2).Value = "Bunner"
$ws.Cells.Item(529, 3).Value = "DRIVING IN MARKED LANES"
$c_529_D = $ws.Cells.Item(529, 4)
$c_529_D.NumberFormat = "@"
$c_529_D.Value = "4511.33"
$c_529_D.Style = "Normal"
$ws.Cells.Item(529, 5).Value = "MM"
$ws.Cells.Item(529, 6).Value = "Guilty"
$ws.Cells.Item(529, 7).Value = "Guilty"
$c_529_H = $ws.Cells.Item(529, 8)
$c_529_H.NumberFormat = "@"
$c_529_H.Value = "$ 0"
$c_529_H.Style = "Normal"
$c_529_I = $ws.Cells.Item(529, 9)
$c_529_I.NumberFormat = "@"
$c_529_I.Value = "$ 0"
$c_529_I.Style = "Normal"

# Row 530
$ws.Cells.Item(530, 1).Value = "21TRC08418"
$ws.Cells.Item(530, 2).Value = "Bunner"
$ws.Cells.Item(530, 3).Value = "TURN AND STOP SIGNALS"
$ws.Cells.Item(530, 4).Value = "No Data"
$ws.Cells.Item(530, 5).Value = "MM"
$ws.Cells.Item(530, 6).Value = "Guilty"
$ws.Cells.Item(530, 7).Value = "Guilty"
$c_530_H = $ws.Cells.Item(530, 8)
$c_530_H.NumberFormat = "@"
$c_530_H.Value = "$ 0"
$c_530_H.Style = "Normal"
$c_530_I = $ws.Cells.Item(530, 9)
$c_530_I.NumberFormat = "@"
$c_530_I.Value = "$ 0"
$c_530_I.Style = "Normal"

# Row 531
$ws.Cells.Item(531, 1).Value = "21TRC08418"
$ws.Cells.Item(531, 2).Value = "Bunner"
$ws.Cells.Item(531, 3).Value = "OVI ALCOHOL / DRUGS 1ST"
$ws.Cells.Item(531, 4).Value = "4511.19A1A*"
$ws.Cells.Item(531, 5).Value = "M1"
$ws.Cells.Item(531, 6).Value = "Guilty"
$ws.Cells.Item(531, 7).Value = "Guilty"
$c_531_H = $ws.Cells.Item(531, 8)
$c_531_H.NumberFormat = "@"
$c_531_H.Value = "$ 0"
$c_531_H.Style = "Normal"
$c_531_I = $ws.Cells.Item(531, 9)
$c_531_I.NumberFormat = "@"
$c_531_I.Value = "$ 0"
$c_531_I.Style = "Normal"

# Row 532
$ws.Cells.Item(532, 1).Value = "21TRC08418"
$ws.Cells.Item(532, 2).Value = "Bunner"
$ws.Cells.Item(532, 3).Value = "DRIVING IN MARKED LANES"
$c_532_D = $ws.Cells.Item(532, 4)
$c_532_D.NumberFormat = "@"
$c_532_D.Value = "4511.33"
$c_532_D.Style = "Normal"
$ws.Cells.Item(532, 5).Value = "MM"
$ws.Cells.Item(532, 6).Value = "Guilty"
$ws.Cells.Item(532, 7).Value = "Guilty"
$c_532_H = $ws.Cells.Item(532, 8)
$c_532_H.NumberFormat = "@"
$c_532_H.Value = "$ 0"
$c_532_H.Style = "Normal"
$c_532_I = $ws.Cells.Item(532, 9)
$c_532_I.NumberFormat = "@"
$c_532_I.Value = "$ 0"
$c_532_I.Style = "Normal"
$ws.Cells.Item(532, 10).Value = "None"
$ws.Cells.Item(532, 11).Value = "None"

# Row 533
$ws.Cells.Item(533, 1).Value = "21TRC08418"
$ws.Cells.Item(533, 2).Value = "Bunner"
$ws.Cells.Item(533, 3).Value = "TURN AND STOP SIGNALS"
$ws.Cells.Item(533, 4).Value = "No Data"
$ws.Cells.Item(533, 5).Value = "MM"
$ws.Cells.Item(533, 6).Value = "Guilty"
$ws.Cells.Item(533, 7).Value = "Guilty"
$c_533_H = $ws.Cells.Item(533, 8)
$c_533_H.NumberFormat = "@"
$c_533_H.Value = "$ 0"
$c_533_H.Style = "Normal"
$c_533_I = $ws.Cells.Item(533, 9)
$c_533_I.NumberFormat = "@"
$c_533_I.Value = "$ 0"
$c_533_I.Style = "Normal"
$ws.Cells.Item(533, 10).Value = "None"
$ws.Cells.Item(533, 11).Value = "None"

# Row 534
$ws.Cells.Item(534, 1).Value = "21TRC08418"
$ws.Cells.Item(534, 2).Value = "Bunner"
$ws.Cells.Item(534, 3).Value = "OVI ALCOHOL / DRUGS 1ST"
$ws.Cells.Item(534, 4).Value = "4511.19A1A*"
$ws.Cells.Item(534, 5).Value = "M1"
$ws.Cells.Item(534, 6).Value = "Guilty"
$ws.Cells.Item(534, 7).Value = "Guilty"
$c_534_H = $ws.Cells.Item(534, 8)
$c_534_H.NumberFormat = "@"
$c_534_H.Value = "$ 0"
$c_534_H.Style = "Normal"
$c_534_I = $ws.Cells.Item(534, 9)
$c_534_I.NumberFormat = "@"
$c_534_I.Value = "$ 0"
$c_534_I.Style = "Normal"
$ws.Cells.Item(534, 10).Value = "None"
$ws.Cells.Item(534, 11).Value = "None"

# Row 535
$ws.Cells.Item(535, 1).Value = "21CRB01291"
$ws.Cells.Item(535, 2).Value = "Hemmeter"
$ws.Cells.Item(535, 3).Value = "PERMISSION REQ'D TO USE LICENSED DOCK"
$ws.Cells.Item(535, 4).Value = "1501:46-12-04"
$ws.Cells.Item(535, 5).Value = "MM"
$ws.Cells.Item(535, 6).Value = "Guilty"
$ws.Cells.Item(535, 7).Value = "Guilty"
$c_535_H = $ws.Cells.Item(535, 8)
$c_535_H.NumberFormat = "@"
$c_535_H.Value = "$ 0"
$c_535_H.Style = "Normal"
$c_535_I = $ws.Cells.Item(535, 9)
$c_535_I.NumberFormat = "@"
$c_535_I.Value = "$ 0"
$c_535_I.Style = "Normal"
$ws.Cells.Item(535, 10).Value = "None"
$ws.Cells.Item(535, 11).Value = "None"

# Row 536
$ws.Cells.Item(536, 1).Value = "21CRB01291"
$ws.Cells.Item(536, 2).Value = "Hemmeter"
$ws.Cells.Item(536, 3).Value = "PERMISSION REQ'D TO USE LICENSED DOCK"
$ws.Cells.Item(536, 4).Value = "1501:46-12-04"
$ws.Cells.Item(536, 5).Value = "MM"
$ws.Cells.Item(536, 6).Value = "Guilty"
$ws.Cells.Item(536, 7).Value = "Guilty"
$c_536_H = $ws.Cells.Item(536, 8)
$c_536_H.NumberFormat = "@"
$c_536_H.Value = "$ 0"
$c_536_H.Style = "Normal"
$c_536_I = $ws.Cells.Item(536, 9)
$c_536_I.NumberFormat = "@"
$c_536_I.Value = "$ 0"
$c_536_I.Style = "Normal"
$ws.Cells.Item(536, 10).Value = "None"
$ws.Cells.Item(536, 11).Value = "None"

# Row 537
$ws.Cells.Item(537, 1).Value = "21CRB01291"
$ws.Cells.Item(537, 2).Value = "Hemmeter"
$ws.Cells.Item(537, 3).Value = "PERMISSION REQ'D TO USE LICENSED DOCK"
$ws.Cells.Item(537, 4).Value = "1501:46-12-04"
$ws.Cells.Item(537, 5).Value = "MM"
$ws.Cells.Item(537, 6).Value = "Guilty"
$ws.Cells.Item(537, 7).Value = "Guilty"
$c_537_H = $ws.Cells.Item(537, 8)
$c_537_H.NumberFormat = "@"
$c_537_H.Value = "$ 0"
$c_537_H.Style = "Normal"
$c_537_I = $ws.Cells.Item(537, 9)
$c_537_I.NumberFormat = "@"
$c_537_I.Value = "$ 0"
$c_537_I.Style = "Normal"
$ws.Cells.Item(537, 10).Value = "None"
$ws.Cells.Item(537, 11).Value = "None"

# Row 538
$ws.Cells.Item(538, 1).Value = "21CRB01291"
$ws.Cells.Item(538, 2).Value = "Hemmeter"
$ws.Cells.Item(538, 3).Value = "PERMISSION REQ'D TO USE LICENSED DOCK"
$ws.Cells.Item(538, 4).Value = "1501:46-12-04"
$ws.Cells.Item(538, 5).Value = "MM"
$ws.Cells.Item(538, 6).Value = "Guilty"
$ws.Cells.Item(538, 7).Value = "Guilty"
$c_538_H = $ws.Cells.Item(538, 8)
$c_538_H.NumberFormat = "@"
$c_538_H.Value = "$ 0"
$c_538_H.Style = "Normal"
$c_538_I = $ws.Cells.Item(538, 9)
$c_538_I.NumberFormat = "@"
$c_538_I.Value = "$ 0"
$c_538_I.Style = "Normal"

# Row 539
$ws.Cells.Item(539, 1).Value = "21CRB01291"
$ws.Cells.Item(539, 2).Value = "Hemmeter"
$ws.Cells.Item(539, 3).Value = "PERMISSION REQ'D TO USE LICENSED DOCK"
$ws.Cells.Item(539, 4).Value = "1501:46-12-04"
$ws.Cells.Item(539, 5).Value = "MM"
$ws.Cells.Item(539, 6).Value = "Guilty"
$ws.Cells.Item(539, 7).Value = "Guilty"
$c_539_H = $ws.Cells.Item(539, 8)
$c_539_H.NumberFormat = "@"
$c_539_H.Value = "$ 0"
$c_539_H.Style = "Normal"
$c_539_I = $ws.Cells.Item(539, 9)
$c_539_I.NumberFormat = "@"
$c_539_I.Value = "$ 0"
$c_539_I.Style = "Normal"

# Row 540
$ws.Cells.Item(540, 1).Value = "21CRB01437"
$ws.Cells.Item(540, 2).Value = "Hemmeter"
$ws.Cells.Item(540, 3).Value = "POSSESSION OF MARIHUANA"
$ws.Cells.Item(540, 4).Value = "2925.11C3"
$ws.Cells.Item(540, 5).Value = "MM"
$ws.Cells.Item(540, 6).Value = "Guilty"
$ws.Cells.Item(540, 7).Value = "Guilty"
$c_540_H = $ws.Cells.Item(540, 8)
$c_540_H.NumberFormat = "@"
$c_540_H.Value = "$ 0"
$c_540_H.Style = "Normal"
$c_540_I = $ws.Cells.Item(540, 9)
$c_540_I.NumberFormat = "@"
$c_540_I.Value = "$ 0"
$c_540_I.Style = "Normal"

# Row 541
$ws.Cells.Item(541, 1).Value = "21CRB01437"
$ws.Cells.Item(541, 2).Value = "Hemmeter"
$ws.Cells.Item(541, 3).Value = "POSSESSION OF MARIHUANA - AMENDED to Assured Clear Distrance Ahead"
$ws.Cells.Item(541, 4).Value = "2925.11C3"
$ws.Cells.Item(541, 5).Value = "MM"
$ws.Cells.Item(541, 6).Value = "Guilty"
$ws.Cells.Item(541, 7).Value = "Guilty"
$c_541_H = $ws.Cells.Item(541, 8)
$c_541_H.NumberFormat = "@"
$c_541_H.Value = "$ 0"
$c_541_H.Style = "Normal"
$c_541_I = $ws.Cells.Item(541, 9)
$c_541_I.NumberFormat = "@"
$c_541_I.Value = "$ 0"
$c_541_I.Style = "Normal"

# Row 542
$ws.Cells.Item(542, 1).Value = "21TRD09246"
$ws.Cells.Item(542, 2).Value = "Bunner"
$ws.Cells.Item(542, 3).Value = "1ST SPEED IN 1 YR >70MPH"
$ws.Cells.Item(542, 4).Value = "4511.21D4"
$ws.Cells.Item(542, 5).Value = "No Data"
$ws.Cells.Item(542, 6).Value = "Guilty"
$ws.Cells.Item(542, 7).Value = "Guilty"
$c_542_H = $ws.Cells.Item(542, 8)
$c_542_H.NumberFormat = "@"
$c_542_H.Value = "$ 0"
$c_542_H.Style = "Normal"
$c_542_I = $ws.Cells.Item(542, 9)
$c_542_I.NumberFormat = "@"
$c_542_I.Value = "$ 0"
$c_542_I.Style = "Normal"
$ws.Cells.Item(542, 10).Value = "None"
$ws.Cells.Item(542, 11).Value = "None"

# Row 543
$ws.Cells.Item(543, 1).Value = "21TRD09246"
$ws.Cells.Item(543, 2).Value = "Bunner"
$ws.Cells.Item(543, 3).Value = "1ST SPEED IN 1 YR >70MPH"
$ws.Cells.Item(543, 4).Value = "4511.21D4"
$ws.Cells.Item(543, 5).Value = "No Data"
$ws.Cells.Item(543, 6).Value = "Guilty"
$ws.Cells.Item(543, 7).Value = "Guilty"
$c_543_H = $ws.Cells.Item(543, 8)
$c_543_H.NumberFormat = "@"
$c_543_H.Value = "$ 0"
$c_543_H.Style = "Normal"
$c_543_I = $ws.Cells.Item(543, 9)
$c_543_I.NumberFormat = "@"
$c_543_I.Value = "$ 0"
$c_543_I.Style = "Normal"
$ws.Cells.Item(543, 10).Value = "None"
$ws.Cells.Item(543, 11).Value = "None"
